$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("web 121")

# Row 3: G3 and I3 go from 0 to 9.5
$ws.Range("G3").Value = 9.5
$ws.Range("I3").Value = 9.5

# Row 14 (student #13) entirely cleared out (values/formulas removed, formatting kept)
$ws.Range("A14:K14").ClearContents()

# Final-scores block (rows 22-33): "Project" column D updated for several students
$ws.Range("D22").Value = 45
$ws.Range("D25").Value = 50
$ws.Range("D26").Value = 56
$ws.Range("D27").Value = 50
$ws.Range("D33").Value = 54

# Row 34 (student #13 final score row) entirely cleared out
$ws.Range("A34:F34").ClearContents()

# Restore the selection to where the editor left off
[void]$ws.Range("D23").Select()
